$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44449, 3, 3, 36.46086533787069),
    @(44450, 1, 4, 48.61448711716091),
    @(44451, 5, 9, 109.3825960136121),
    @(44452, 0, 9, 109.3825960136121),
    @(44453, 0, 9, 109.3825960136121),
    @(44454, 0, 9, 109.3825960136121),
    @(44455, 1, 10, 121.5362177929023),
    @(44456, 1, 8, 97.22897423432183),
    @(44457, 0, 7, 85.07535245503161),
    @(44458, 0, 2, 24.30724355858046),
    @(44459, 1, 3, 36.46086533787069)
)

$startRow = 375
$lastRow = $startRow - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy formatting (style, number format) from the row above, like Excel
    # autofilling a date-formatted column when new data is appended.
    $ws.Cells.Item($lastRow, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
